# Team1_May_06_2009_Final_Presentation.pptx edit
#
# 1. Bump the cached "datetimeFigureOut" date placeholder text from
#    5/1/2009 to 5/4/2009 on the slide master and every slide layout.
# 2. Fix the typo "Collaboracion" -> "Collaboracíon" in the title of the
#    last slide ("Viva la Collaboracion!").

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, $label) {
    $count = $shapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $sh = $shapes.Item($i)
        $phType = -1
        try { $phType = $sh.PlaceholderFormat.Type } catch {}
        if ($phType -eq 16) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "5/1/2009") {
                $tr.Text = "5/4/2009"
            }
        }
    }
}

# -- Slide master --
$master = $p.SlideMaster
$masterShapes = $master.Shapes
Update-DatePlaceholder $masterShapes "Master"

# -- Every slide layout --
$layouts = $master.CustomLayouts
$layoutCount = $layouts.Count
for ($li = 1; $li -le $layoutCount; $li++) {
    $lay = $layouts.Item($li)
    $layShapes = $lay.Shapes
    $label = "Layout " + $li
    Update-DatePlaceholder $layShapes $label
}

# -- Fix "Collaboracion" -> "Collaboracíon" typo on the last slide's title --
$slideCount = $p.Slides.Count
$lastSlide = $p.Slides.Item($slideCount)
$shapes = $lastSlide.Shapes
$shapeCount = $shapes.Count
for ($i = 1; $i -le $shapeCount; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $fullText = $tr.Text
        $idx = $fullText.IndexOf("Collaboracion")
        if ($idx -ge 0) {
            $word = $tr.Characters($idx + 1, 13)
            $word.Text = "Collaboracíon"
        }
    }
}
